$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D19").Value = "2016-03-02 15:39:00"
$wsZh.Range("D20").Value = "2016-03-02 15:39:00"
$wsZh.Range("G19").Value = "2016-03-02 15:39:47"
$wsZh.Range("G20").Value = "2016-03-02 15:39:47"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D19").Value = "2016-03-02 15:39:11"
$wsDe.Range("D20").Value = "2016-03-02 15:39:11"
$wsDe.Range("G19").Value = "2016-03-02 15:40:07"
$wsDe.Range("G20").Value = "2016-03-02 15:40:07"
